$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '44.181.56'
$ws.Range('E2').Value = '  +2.41%  '
$ws.Range('D3').Value = '2.427.49'
$ws.Range('E3').Value = '  +2.10%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '308.27'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.71%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '100.93'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.24%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.514'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.500'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.02%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.36'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.29%  '
$ws.Range('E11').Value = '  +1.70%  '
$ws.Range('E12').Value = '  +2.56%  '
$ws.Range('E13').Value = '  +1.76%  '
$ws.Range('E14').Value = '  +2.21%  '
$ws.Range('D15').Value = '2.804.40'
$ws.Range('E15').Value = '  +1.91%  '
$ws.Range('D16').Value = '2.424.14'
$ws.Range('E16').Value = '  +1.34%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.835'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +3.50%  '
$ws.Range('D18').Value = '44.109.90'
$ws.Range('E18').Value = '  +2.23%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.33'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.38%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.43'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.72%  '
$ws.Range('D21').Value = '0.0₃0906'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '68.58'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.01%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '240.52'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.28%  '
$ws.Range('E24').Value = '  +2.39%  '
$ws.Range('E25').Value = '  +1.30%  '
$ws.Range('E26').Value = '  -0.02%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '25.24'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.55%  '
$ws.Range('E28').Value = '  -1.10%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.62'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +5.36%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '33.24'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +5.56%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.118'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +15.04%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '18.65'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +8.90%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.18'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.65%  '
$ws.Range('E34').Value = '  +0.00%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0765'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +4.10%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.91'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.62%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.49'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +4.80%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '129.11'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +24.53%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.92'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +4.56%  '
$ws.Range('E40').Value = '  -0.12%  '
$ws.Range('E41').Value = '  +0.27%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '21.24'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -5.24%  '
$ws.Range('E43').Value = '  +3.04%  '
$ws.Range('D44').Value = '1.957.36'
$ws.Range('E44').Value = '  +0.12%  '
$ws.Range('E45').Value = '  +1.69%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.88'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +4.84%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.42'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.41%  '
$ws.Range('E48').Value = '  +9.09%  '
$ws.Range('B49').Value = 'RocketPoolETH'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D49').Value = '2.668.86'
$ws.Range('E49').Value = '  +2.28%  '
$ws.Range('B50').Value = 'MultiversX'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '53.49'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.45%  '
$ws.Range('B51').Value = 'BitcoinSV'
$ws.Range('C51').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '73.49'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.40%  '
